$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Crossword Puzzle Comp award: should be $150, not $100
$ws.Range("E17").Value = 150

# Reflect the last active selection recorded in the saved workbook
$ws.Range("E18").Select()
